$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7 - everything from the old row 7 onward
# (old header row 8, data rows 9-31, trailing formatting rows, etc.) shifts
# down by one row.
$ws.Rows("7:7").Insert()

# Reset inherited formatting (the freshly inserted row copies the format of
# the row above it) before applying the new "intended use" look.
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Style = "Normal"

# Populate the new row with the intended-use label/value pair.
$ws.Range("B7").Value = "intendedUse"
$ws.Range("C7").Value = "Epi-validated outbreak"

# Style: label is bold black, value is plain black left-aligned.
$ws.Range("B7").Font.Bold = $true
$ws.Range("B7").Font.Color = 0
$ws.Range("C7").Font.Color = 0
$ws.Range("C7").HorizontalAlignment = -4131

# The named range "tmpt." pointed at the data block $D$9:$D$31; since the
# data moved down one row with the insert, repoint it at $D$10:$D$32.
$ws.Names.Item(1).RefersTo = "='Salmonella_enterica_1203NYJAP-1'!`$D`$10:`$D`$32"

# Restore the selection to the cell that was being edited.
[void]$ws.Range("E10").Select()
